# Adds two new data rows (4 and 5) to the Artfynd sheet, matching the
# records that were appended to the source OOXML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value as genuine text (quotePrefix), even when the
# text looks like a number or a date (e.g. "1", "2025", "2024-05-03").
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
}

# =======================================================================
# Row 4
# =======================================================================
$ws.Range("A4").Value = 131124849
$ws.Range("B4").Value = 8261
Set-TextValue $ws.Range("D4") "NT"
$ws.Range("E4").Value = 106456
Set-TextValue $ws.Range("F4") "Granvivel"
Set-TextValue $ws.Range("G4") "Pissodes harcyniae"
Set-TextValue $ws.Range("H4") "(Herbst, 1795)"
Set-TextValue $ws.Range("I4") "1"
Set-TextValue $ws.Range("J4") "ex."
Set-TextValue $ws.Range("K4") "imago/adult"
Set-TextValue $ws.Range("N4") "fönsterfälla"
Set-TextValue $ws.Range("P4") "Fäbodsberget, N, Upl"
$ws.Range("Q4").Value = 619246
$ws.Range("R4").Value = 6662712
$ws.Range("S4").Value = 25
Set-TextValue $ws.Range("T4") "Uppsala"
Set-TextValue $ws.Range("U4") "Heby"
Set-TextValue $ws.Range("V4") "Uppland"
Set-TextValue $ws.Range("W4") "Harbo"
Set-TextValue $ws.Range("Y4") "2024-05-03"
Set-TextValue $ws.Range("AA4") "2024-06-04"
Set-TextValue $ws.Range("AC4") "IBL-fälla"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
Set-TextValue $ws.Range("AQ4") "Hans-Erik Wanntorp"
Set-TextValue $ws.Range("AS4") "Hans-Erik Wanntorp"
Set-TextValue $ws.Range("AT4") "2025"
Set-TextValue $ws.Range("AW4") "Hans-Erik Wanntorp"
Set-TextValue $ws.Range("AX4") "Joachim Strengbom"

# =======================================================================
# Row 5
# =======================================================================
$ws.Range("A5").Value = 131124909
$ws.Range("B5").Value = 6282
Set-TextValue $ws.Range("D5") "NT"
$ws.Range("E5").Value = 100524
Set-TextValue $ws.Range("F5") "Skrovlig flatbagge"
Set-TextValue $ws.Range("G5") "Calitys scabra"
Set-TextValue $ws.Range("H5") "(Thunberg, 1784)"
Set-TextValue $ws.Range("I5") "1"
Set-TextValue $ws.Range("J5") "ex."
Set-TextValue $ws.Range("K5") "imago/adult"
Set-TextValue $ws.Range("N5") "fönsterfälla"
Set-TextValue $ws.Range("P5") "Fäbodsberget, N, Upl"
$ws.Range("Q5").Value = 619246
$ws.Range("R5").Value = 6662712
$ws.Range("S5").Value = 25
Set-TextValue $ws.Range("T5") "Uppsala"
Set-TextValue $ws.Range("U5") "Heby"
Set-TextValue $ws.Range("V5") "Uppland"
Set-TextValue $ws.Range("W5") "Harbo"
Set-TextValue $ws.Range("Y5") "2024-05-03"
Set-TextValue $ws.Range("AA5") "2024-06-04"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
Set-TextValue $ws.Range("AI5") "produktionsskog"
Set-TextValue $ws.Range("AS5") "Hans-Erik Wanntorp"
Set-TextValue $ws.Range("AW5") "Hans-Erik Wanntorp"
Set-TextValue $ws.Range("AX5") "Joachim Strengbom"

Write-Host "Rows 4 and 5 written."
